$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F4").Value = "application instructions"
$ws.Range("F6").Value = "env warning - water || off target movement"
$ws.Range("F8").Value = "154_pesticide_storage"
$ws.Range("F13").Value = "use restrictions"
$ws.Range("F14").Value = "application instructions"
$ws.Range("F16").Value = "safety procedures"
$ws.Range("F19").Value = "chemigation"
$ws.Range("F20").Value = "irrigation"
$ws.Range("F22").Value = "mixing"
$ws.Range("F23").Value = "mixing"
